$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2800.651
$ws.Cells.Item(138, 9).Value = 1184.4634
$ws.Cells.Item(138, 10).Value = 5812.636
$ws.Cells.Item(138, 11).Value = 3553.3902
$ws.Cells.Item(138, 12).Value = 17437.908
$ws.Cells.Item(138, 13).Value = 1586.6098
$ws.Cells.Item(138, 14).Value = -27717.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 6263.52
$ws.Cells.Item(74, 9).Value = 2922.7144
$ws.Cells.Item(74, 10).Value = 10515.454
$ws.Cells.Item(74, 11).Value = 2922.7144
$ws.Cells.Item(74, 12).Value = 10515.454
$ws.Cells.Item(74, 13).Value = -2048.7144
$ws.Cells.Item(74, 14).Value = -12263.454

$ws.Cells.Item(77, 8).Value = 6263.52
$ws.Cells.Item(77, 9).Value = 2922.7144
$ws.Cells.Item(77, 10).Value = 10515.454
$ws.Cells.Item(77, 11).Value = 14613.572
$ws.Cells.Item(77, 12).Value = 52577.27
$ws.Cells.Item(77, 13).Value = -10245.572
$ws.Cells.Item(77, 14).Value = -61313.27

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(18, 8).Value = 13840
$ws.Cells.Item(18, 10).Value = 13840
$ws.Cells.Item(18, 12).Value = 13840
$ws.Cells.Item(18, 14).Value = -14898

$ws.Cells.Item(86, 8).Value = 1634.0731
$ws.Cells.Item(86, 9).Value = 1628.9474
$ws.Cells.Item(86, 10).Value = 1699
$ws.Cells.Item(86, 11).Value = 1628.9474
$ws.Cells.Item(86, 12).Value = 1699
$ws.Cells.Item(86, 13).Value = -505.9474
$ws.Cells.Item(86, 14).Value = -3945

$ws.Cells.Item(89, 8).Value = 1634.0731
$ws.Cells.Item(89, 9).Value = 1628.9474
$ws.Cells.Item(89, 10).Value = 1699
$ws.Cells.Item(89, 11).Value = 8144.737
$ws.Cells.Item(89, 12).Value = 8495
$ws.Cells.Item(89, 13).Value = -2528.737
$ws.Cells.Item(89, 14).Value = -19727

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2103.532
$ws.Cells.Item(31, 9).Value = 1383.1143
$ws.Cells.Item(31, 10).Value = 4204.75
$ws.Cells.Item(31, 11).Value = 1383.1143
$ws.Cells.Item(31, 12).Value = 4204.75
$ws.Cells.Item(31, 13).Value = -1088.1143
$ws.Cells.Item(31, 14).Value = -4794.75

$ws.Cells.Item(34, 8).Value = 2103.532
$ws.Cells.Item(34, 9).Value = 1383.1143
$ws.Cells.Item(34, 10).Value = 4204.75
$ws.Cells.Item(34, 11).Value = 1383.1143
$ws.Cells.Item(34, 12).Value = 4204.75
$ws.Cells.Item(34, 13).Value = -1181.1143
$ws.Cells.Item(34, 14).Value = -4608.75

$ws.Cells.Item(94, 8).Value = 1508.7142
$ws.Cells.Item(94, 9).Value = 1299
$ws.Cells.Item(94, 10).Value = 1543.6666
$ws.Cells.Item(94, 11).Value = 1299
$ws.Cells.Item(94, 12).Value = 1543.6666
$ws.Cells.Item(94, 13).Value = -848
$ws.Cells.Item(94, 14).Value = -2445.6666

$ws.Cells.Item(99, 8).Value = 3568.5
$ws.Cells.Item(99, 9).Value = 3046.889
$ws.Cells.Item(99, 10).Value = 5133.3335
$ws.Cells.Item(99, 11).Value = 3046.889
$ws.Cells.Item(99, 12).Value = 5133.3335
$ws.Cells.Item(99, 13).Value = -1548.889
$ws.Cells.Item(99, 14).Value = -8129.3335

$ws.Cells.Item(126, 8).Value = 3568.5
$ws.Cells.Item(126, 9).Value = 3046.889
$ws.Cells.Item(126, 10).Value = 5133.3335
$ws.Cells.Item(126, 11).Value = 9140.667000000001
$ws.Cells.Item(126, 12).Value = 15400.0005
$ws.Cells.Item(126, 13).Value = -6670.667000000001
$ws.Cells.Item(126, 14).Value = -20340.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1176.1719
$ws.Cells.Item(131, 9).Value = 1778
$ws.Cells.Item(131, 10).Value = 1064.7222
$ws.Cells.Item(131, 11).Value = 5334
$ws.Cells.Item(131, 12).Value = 3194.1666
$ws.Cells.Item(131, 13).Value = -294
$ws.Cells.Item(131, 14).Value = -13274.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(23, 8).Value = 7365.4287
$ws.Cells.Item(23, 9).Value = 3751
$ws.Cells.Item(23, 10).Value = 8811.200000000001
$ws.Cells.Item(23, 11).Value = 3751
$ws.Cells.Item(23, 12).Value = 8811.200000000001
$ws.Cells.Item(23, 13).Value = -3528
$ws.Cells.Item(23, 14).Value = -9257.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 4935.1665
$ws.Cells.Item(9, 9).Value = 366.66666
$ws.Cells.Item(9, 10).Value = 9503.666999999999
$ws.Cells.Item(9, 11).Value = 366.66666
$ws.Cells.Item(9, 12).Value = 9503.666999999999
$ws.Cells.Item(9, 13).Value = -142.66666
$ws.Cells.Item(9, 14).Value = -9951.666999999999

$ws.Cells.Item(16, 8).Value = 1001.4091
$ws.Cells.Item(16, 9).Value = 627.7222
$ws.Cells.Item(16, 10).Value = 2683
$ws.Cells.Item(16, 11).Value = 627.7222
$ws.Cells.Item(16, 12).Value = 2683
$ws.Cells.Item(16, 13).Value = -457.7222
$ws.Cells.Item(16, 14).Value = -3023

$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).ClearContents()

$ws.Cells.Item(22, 8).Value = 1142.8462
$ws.Cells.Item(22, 9).Value = 1485
$ws.Cells.Item(22, 10).Value = 1040.2
$ws.Cells.Item(22, 11).Value = 1485
$ws.Cells.Item(22, 12).Value = 1040.2
$ws.Cells.Item(22, 13).Value = -1190
$ws.Cells.Item(22, 14).Value = -1630.2

$ws.Cells.Item(23, 8).Value = 19800
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).ClearContents()

$ws.Cells.Item(24, 8).Value = 4000
$ws.Cells.Item(24, 10).Value = 4000
$ws.Cells.Item(24, 12).Value = 4000
$ws.Cells.Item(24, 14).Value = -4686

$ws.Cells.Item(25, 8).Value = 60004
$ws.Cells.Item(25, 10).Value = 60004
$ws.Cells.Item(25, 12).Value = 60004
$ws.Cells.Item(25, 14).Value = -60464

$ws.Cells.Item(27, 8).Value = 1142.8462
$ws.Cells.Item(27, 9).Value = 1485
$ws.Cells.Item(27, 10).Value = 1040.2
$ws.Cells.Item(27, 11).Value = 1485
$ws.Cells.Item(27, 12).Value = 1040.2
$ws.Cells.Item(27, 13).Value = -1378
$ws.Cells.Item(27, 14).Value = -1254.2

$ws.Cells.Item(32, 8).Value = 12900
$ws.Cells.Item(32, 10).Value = 24800
$ws.Cells.Item(32, 12).Value = 24800
$ws.Cells.Item(32, 14).Value = -25434

$ws.Cells.Item(46, 8).Value = 1778.5
$ws.Cells.Item(46, 9).Value = 1996.4
$ws.Cells.Item(46, 10).Value = 1622.8572
$ws.Cells.Item(46, 11).Value = 1996.4
$ws.Cells.Item(46, 12).Value = 1622.8572
$ws.Cells.Item(46, 13).Value = -1808.4
$ws.Cells.Item(46, 14).Value = -1998.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 7362.6665
$ws.Cells.Item(39, 9).Value = 5044
$ws.Cells.Item(39, 10).Value = 12000
$ws.Cells.Item(39, 11).Value = 5044
$ws.Cells.Item(39, 12).Value = 12000
$ws.Cells.Item(39, 13).Value = -4631
$ws.Cells.Item(39, 14).Value = -12826

$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 14).ClearContents()

$ws.Cells.Item(43, 8).Value = 19576.666
$ws.Cells.Item(43, 9).Value = 1500
$ws.Cells.Item(43, 10).Value = 28615
$ws.Cells.Item(43, 11).Value = 1500
$ws.Cells.Item(43, 12).Value = 28615
$ws.Cells.Item(43, 13).Value = -1351
$ws.Cells.Item(43, 14).Value = -28913

$ws.Cells.Item(96, 8).Value = 1500
$ws.Cells.Item(96, 10).Value = 1500
$ws.Cells.Item(96, 12).Value = 1500
$ws.Cells.Item(96, 14).Value = -4246
